$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# Insert a new row at position 5; this shifts existing rows 5-16 down to 6-17
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).ClearFormats()

# Populate the new row 5 with the new trade data
$ws.Range("A5").Value = 46062
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 301.2
$ws.Range("F5").Value = 3027.1
$ws.Range("G5").Value = "CN#252611665409"
$ws.Range("I5").Value = 15.1
$ws.Range("J5").Formula = '=Index!$C$2'

# Match the date style of the other DATE column cells (A6 has the date number format)
$ws.Range("A6").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
